$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "13-jul" day column (AC) after the existing "12-jul" column (AB)
$ws.Range("AC1").Value = "13-jul"

$ws.Range("AC2").Value  = 0
$ws.Range("AC3").Value  = 11.863420027590076
$ws.Range("AC4").Value  = 18.941985896417894
$ws.Range("AC5").Value  = 24.747402788658988
$ws.Range("AC6").Value  = 0
$ws.Range("AC7").Value  = 10.957853298839492
$ws.Range("AC8").Value  = 10.542004601093442
$ws.Range("AC9").Value  = 23.277091139825085
$ws.Range("AC10").Value = 21.26781403911172
$ws.Range("AC11").Value = 12.67597142544755
$ws.Range("AC12").Value = 0
$ws.Range("AC13").Value = 11.082497098897598
$ws.Range("AC14").Value = 0
$ws.Range("AC15").Value = 0
$ws.Range("AC16").Value = 18.497377502682234
$ws.Range("AC17").Value = 0
$ws.Range("AC18").Value = 0

# Update the selection to match the saved workbook state (cell AC14 active)
[void]$ws.Range("AC14").Select()
